$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.570.76'
$ws.Range('E2').Value = '  -3.56%  '

$ws.Range('D3').Value = '2.377.00'
$ws.Range('E3').Value = '  -4.52%  '

$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').Value = "'509.87"
$ws.Range('E5').Value = '  -4.38%  '

$ws.Range('D6').Value = "'130.08"
$ws.Range('E6').Value = '  -2.43%  '

$ws.Range('D7').Value = "'0.998"
$ws.Range('E7').Value = '  -0.69%  '

$ws.Range('E8').Value = '  -1.96%  '

$ws.Range('D9').Value = '2.401.51'
$ws.Range('E9').Value = '  -4.00%  '

$ws.Range('D10').Value = "'0.0966"
$ws.Range('E10').Value = '  -2.32%  '

$ws.Range('D11').Value = "'0.150"
$ws.Range('E11').Value = '  -1.87%  '

$ws.Range('E12').Value = '  -2.16%  '

$ws.Range('D13').Value = "'4.69"
$ws.Range('E13').Value = '  -9.75%  '

$ws.Range('D14').Value = '2.802.57'
$ws.Range('E14').Value = '  -4.55%  '

$ws.Range('D15').Value = '56.473.66'
$ws.Range('E15').Value = '  -3.69%  '

$ws.Range('D16').Value = "'21.64"
$ws.Range('E16').Value = '  -2.67%  '

$ws.Range('E17').Value = '  -2.64%  '

$ws.Range('D18').Value = '2.398.42'
$ws.Range('E18').Value = '  -4.38%  '

$ws.Range('D19').Value = "'10.25"
$ws.Range('E19').Value = '  -3.06%  '

$ws.Range('D20').Value = "'313.13"
$ws.Range('E20').Value = '  -2.13%  '

$ws.Range('E21').Value = '  -4.02%  '

$ws.Range('D22').Value = "'6.30"
$ws.Range('E22').Value = '  +1.87%  '

$ws.Range('E23').Value = '  -0.06%  '

$ws.Range('D24').Value = "'65.68"
$ws.Range('E24').Value = '  -0.13%  '

$ws.Range('E25').Value = '  -0.30%  '

$ws.Range('D26').Value = '2.494.47'
$ws.Range('E26').Value = '  -5.06%  '

$ws.Range('B27').Value = 'Polygon'
$ws.Range('C27').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D27').Value = "'0.377"
$ws.Range('E27').Value = '  -7.41%  '

$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').Value = "'0.152"
$ws.Range('E28').Value = '  -4.36%  '

$ws.Range('D29').Value = "'7.24"
$ws.Range('E29').Value = '  -2.91%  '

$ws.Range('D30').Value = "'174.96"
$ws.Range('E30').Value = '  +1.42%  '

$ws.Range('E31').Value = '  -2.45%  '

$ws.Range('E32').Value = '  -4.76%  '

$ws.Range('D33').Value = "'6.14"
$ws.Range('E33').Value = '  -2.39%  '

$ws.Range('E34').Value = '  -5.77%  '

$ws.Range('E35').Value = '  -0.18%  '

$ws.Range('D36').Value = "'0.995"
$ws.Range('E36').Value = '  -0.14%  '

$ws.Range('E37').Value = '  -1.99%  '

$ws.Range('E38').Value = '  -1.92%  '

$ws.Range('E39').Value = '  -4.91%  '

$ws.Range('D40').Value = "'35.85"
$ws.Range('E40').Value = '  -1.41%  '

$ws.Range('D41').Value = "'1.43"
$ws.Range('E41').Value = '  -4.66%  '

$ws.Range('D42').Value = "'0.787"
$ws.Range('E42').Value = '  -2.76%  '

$ws.Range('D43').Value = "'134.20"
$ws.Range('E43').Value = '  +2.47%  '

$ws.Range('E44').Value = '  -3.14%  '

$ws.Range('E45').Value = '  -1.97%  '

$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').Value = "'0.573"
$ws.Range('E46').Value = '  -2.65%  '

$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = "'255.13"
$ws.Range('E47').Value = '  -7.00%  '

$ws.Range('D48').Value = "'0.0902"
$ws.Range('E48').Value = '  -3.23%  '

$ws.Range('D49').Value = "'0.0489"
$ws.Range('E49').Value = '  -4.19%  '

$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'16.83"
$ws.Range('E50').Value = '  -4.18%  '

$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').Value = "'0.0207"
$ws.Range('E51').Value = '  -4.65%  '

